$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 102246.19
$ws.Range("C2").Value = 13791.47

$ws.Range("B3").Value = 910258.05
$ws.Range("C3").Value = 46548.14

$ws.Range("B4").Value = 1749720.49
$ws.Range("C4").Value = 22690.39

$ws.Range("B5").Value = 2824609.64
$ws.Range("C5").Value = 21541.02

$ws.Range("B6").Value = 4430583.14
$ws.Range("C6").Value = 26899.54

$ws.Range("B7").Value = 1974408.92
$ws.Range("C7").Value = 14642.46
